# Rework the single "sample" table sheet into a 3-sheet "methods & tables"
# workbook: Datasets / Controls / Outcomes.

$wb = $excel.ActiveWorkbook

# --- Rename & repurpose the original sheet into "Datasets" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Datasets"

# Only wipe cell contents (B/C of rows 2:3 and all of row 4:5) - keep the
# D2/D3 formatting (vertical-center style) already on those cells.
$ws1.Range("A1:D5").ClearContents()

$ws1.Range("A1").Value = "Dataset"
$ws1.Range("B1").Value = "Description"
$ws1.Range("C1").Value = "Example variables"

$ws1.Range("A2").Value = "Opal usage"
$ws1.Range("A3").Value = "Visit data"

$ws1.Rows.Item(2).RowHeight = 16
$ws1.Rows.Item(3).RowHeight = 16

# ColumnWidth (chars) -> stored <col width> has a constant +5/6 padding
# baked in by the writer, so back it out to land on the exact target widths.
$ws1.Columns.Item(1).ColumnWidth = 9.5 - (5/6)
$ws1.Columns.Item(2).ColumnWidth = 10 - (5/6)
$ws1.Columns.Item(3).ColumnWidth = 15 - (5/6)

$ws1.Application.ActiveWindow.Zoom = 180
$ws1.Range("C1").Select() | Out-Null

# --- Add "Controls" sheet right after "Datasets" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Controls"
$ws2.Range("A1").Value = "Variable"
$ws2.Range("B1").Value = "Description/preprocessing notes"
$ws2.Columns.Item(2).ColumnWidth = 31.5 - (5/6)
$ws2.Application.ActiveWindow.Zoom = 150
$ws2.Range("B3").Select() | Out-Null

# --- Add "Outcomes" sheet right after "Controls" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Outcomes"
$ws3.Range("A1").Value = "Variable"
$ws3.Range("B1").Value = "Description/preprocessing notes"
$ws3.Columns.Item(2).ColumnWidth = 26 - (5/6)
$ws3.Application.ActiveWindow.Zoom = 160
$ws3.Range("B2").Select() | Out-Null

# --- Ensure sheet order is Datasets, Controls, Outcomes & Datasets is active ---
$ws1.Select() | Out-Null
